# Update selection on Sheet5
$wb = $excel.ActiveWorkbook
$ws5 = $wb.Worksheets.Item("Sheet5")
$ws5.Activate()
$ws5.Range("K15").Select()

# Insert a new worksheet before Sheet6, name it Sheet19
$ws6 = $wb.Worksheets.Item("Sheet6")
$newSheet = $wb.Worksheets.Add($ws6)
$newSheet.Name = "Sheet19"

# Header row
$newSheet.Range("A1").Value = "source"
$newSheet.Range("B1").Value = "author"
$newSheet.Range("C1").Value = "table_name"
$newSheet.Range("D1").Value = "dimensions"
$newSheet.Range("E1").Value = "dimension_levels_text"
$newSheet.Range("F1").Value = "dimension_levels_code"
$newSheet.Range("G1").Value = "unit"
$newSheet.Range("H1").Value = "interval"
$newSheet.Range("I1").Value = "series_name"
$newSheet.Range("J1").Value = "table_code"
$newSheet.Range("K1").Value = "series_code"

# Row 2
$newSheet.Range("A2").Value = "umar"
$newSheet.Range("B2").Value = "mz"
$newSheet.Range("C2").Value = "dfg"
$newSheet.Range("D2").Value = "dff"
$newSheet.Range("E2").Value = 3
$newSheet.Range("F2").Value = 234
$newSheet.Range("G2").Value = "%"
$newSheet.Range("H2").Value = "M"
$newSheet.Range("I2").Value = "tesx"
$newSheet.Range("J2").Value = "MZ001"
$newSheet.Range("K2").Value = "UMAR--MZ001--234--M"

# Row 3
$newSheet.Range("A3").Value = "umar"
$newSheet.Range("B3").Value = "mz"
$newSheet.Range("C3").Value = "dfg"
$newSheet.Range("D3").Value = "dff"
$newSheet.Range("E3").Value = 2
$newSheet.Range("F3").Value = 1123
$newSheet.Range("G3").Value = "%"
$newSheet.Range("H3").Value = "M"
$newSheet.Range("I3").Value = "sdt"
$newSheet.Range("J3").Value = "MZ001"
$newSheet.Range("K3").Value = "UMAR--MZ001--1123--M"

# Row 4
$newSheet.Range("A4").Value = "umar"
$newSheet.Range("B4").Value = "mz"
$newSheet.Range("C4").Value = "dfg"
$newSheet.Range("D4").Value = "dff"
$newSheet.Range("E4").Value = 4
$newSheet.Range("F4").Value = 1
$newSheet.Range("G4").Value = "%"
$newSheet.Range("H4").Value = "M"
$newSheet.Range("I4").Value = "sdt"
$newSheet.Range("J4").Value = "MZ002"
$newSheet.Range("K4").Value = "UMAR--MZ002--1--M"

# Row 5
$newSheet.Range("A5").Value = "umar"
$newSheet.Range("B5").Value = "mz"
$newSheet.Range("C5").Value = "dfg"
$newSheet.Range("D5").Value = "dff"
$newSheet.Range("E5").Value = 5
$newSheet.Range("F5").Value = 12
$newSheet.Range("G5").Value = "%"
$newSheet.Range("H5").Value = "M"
$newSheet.Range("I5").Value = "sdt"
$newSheet.Range("J5").Value = "MZ002"
$newSheet.Range("K5").Value = "UMAR--MZ002--12--M"

$newSheet.Range("H6").Select()
